$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.997.04"
$ws.Range("E2").Value = "  +0.26%  "

$ws.Range("D3").Value = "1.642.23"
$ws.Range("E3").Value = "  -0.10%  "

$ws.Range("E4").Value = "  -0.59%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.56"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.31%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5088"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.88%  "

$ws.Range("E7").Value = "  -0.28%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2562"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -0.52%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06347"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -0.80%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.52"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -0.61%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07762"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -0.55%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.286"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +0.11%  "

$ws.Range("D13").Value = "1.642.88"
$ws.Range("E13").Value = "  -0.29%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5425"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -0.20%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "64.20"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -0.90%  "

$ws.Range("D16").Value = "0.0₅7704"
$ws.Range("E16").Value = "  -2.11%  "

$ws.Range("D17").Value = "26.041.80"
$ws.Range("E17").Value = "  +0.34%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.001"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.50%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "198.37"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.12%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.423"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.65%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.913"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -0.60%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.035"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.91%  "

$ws.Range("E23").Value = "  -0.26%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.869"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -0.59%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "141.31"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.70%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1190"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +3.89%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.798"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.91%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.60"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -0.76%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.234"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.62%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.04862"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -0.89%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.251"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.29%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.160"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -1.27%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.524"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -0.82%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.366"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -0.43%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.8985"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +0.48%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.582"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -1.04%  "

$ws.Range("D37").Value = "1.140.10"
$ws.Range("E37").Value = "  +0.01%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5447"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -1.88%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01561"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +0.02%  "

$ws.Range("E40").Value = "  -0.59%  "

$ws.Range("E41").Value = "  +8.62%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.526"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -1.32%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8107"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -0.94%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "99.25"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -0.38%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.384"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -5.53%  "

$ws.Range("D46").Value = "1.783.93"
$ws.Range("E46").Value = "  +0.36%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4527"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -0.20%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "54.90"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.95%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.9990"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -0.99%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05055"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -0.68%  "

$ws.Range("E51").Value = "  -0.42%  "
